$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update the Marking row (row 11) right-answer marks value: 3 -> 5
$ws.Range("B11").Value = 5

# Update the Total row (row 12) right-answer marks value: 45 -> 75
$ws.Range("B12").Value = 75

# Update the Correct/Total marks summary text: 41/84 -> 75/140
$ws.Range("E12").Value = "75/140"
